$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 546.7765096666667
$ws.Range("H2").Value = 1640.329529
$ws.Range("I2").Value = 0.6285526459909564
$ws.Range("J2").Value = 0.6285526459909564
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.987076
$ws.Range("N2").Value = 11.961228
$ws.Range("O2").Value = 0.2813308272685638
$ws.Range("P2").Value = 0.2813308272685638
$ws.Range("Q2").Value = 2180.039499055735
$ws.Range("R2").Value = 19620.35549150161
$ws.Range("S2").Value = 0.1768312358784805
$ws.Range("T2").Value = 0.1768312358784805
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 546.7765096666667
$ws.Range("H3").Value = 1640.329529
$ws.Range("I3").Value = 0.6285526459909564
$ws.Range("J3").Value = 0.6285526459909564
$ws.Range("M3").Value = 10.131229
$ws.Range("N3").Value = 30.393687
$ws.Range("O3").Value = 0.7148664925918803
$ws.Range("P3").Value = 0.7148664925918804
$ws.Range("Q3").Value = 5539.518031253714
$ws.Range("R3").Value = 49855.66228128342
$ws.Range("S3").Value = 0.4493312254489008
$ws.Range("T3").Value = 0.4493312254489009
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 546.7765096666667
$ws.Range("H4").Value = 1640.329529
$ws.Range("I4").Value = 0.6285526459909564
$ws.Range("J4").Value = 0.6285526459909564
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05389233333333333
$ws.Range("N4").Value = 0.161677
$ws.Range("O4").Value = 0.00380268013955587
$ws.Range("P4").Value = 0.00380268013955587
$ws.Range("Q4").Value = 29.46706191779255
$ws.Range("R4").Value = 265.203557260133
$ws.Range("S4").Value = 0.002390184663575102
$ws.Range("T4").Value = 0.002390184663575102
$ws.Range("D5").Value = "FAPs"
$ws.Range("I5").Value = 0.1861770314550556
$ws.Range("J5").Value = 0.1861770314550556
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.987076
$ws.Range("N5").Value = 11.961228
$ws.Range("O5").Value = 0.2813308272685638
$ws.Range("P5").Value = 0.2813308272685638
$ws.Range("Q5").Value = 645.726789916024
$ws.Range("R5").Value = 5811.541109244215
$ws.Range("S5").Value = 0.05237733827765622
$ws.Range("T5").Value = 0.05237733827765622
$ws.Range("D6").Value = "MuSCs"
$ws.Range("I6").Value = 0.1861770314550556
$ws.Range("J6").Value = 0.1861770314550556
$ws.Range("M6").Value = 10.131229
$ws.Range("N6").Value = 30.393687
$ws.Range("O6").Value = 0.7148664925918803
$ws.Range("P6").Value = 0.7148664925918804
$ws.Range("Q6").Value = 1640.802929283046
$ws.Range("R6").Value = 14767.22636354741
$ws.Range("S6").Value = 0.1330917214774438
$ws.Range("T6").Value = 0.1330917214774438
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("I7").Value = 0.1861770314550556
$ws.Range("J7").Value = 0.1861770314550556
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.05389233333333333
$ws.Range("N7").Value = 0.161677
$ws.Range("O7").Value = 0.00380268013955587
$ws.Range("P7").Value = 0.00380268013955587
$ws.Range("Q7").Value = 8.728131443799333
$ws.Range("R7").Value = 78.55318299419399
$ws.Range("S7").Value = 0.0007079716999556085
$ws.Range("T7").Value = 0.0007079716999556086
$ws.Range("D8").Value = "FAPs"
$ws.Range("G8").Value = 160.630483
$ws.Range("H8").Value = 481.891449
$ws.Range("I8").Value = 0.1846544489960017
$ws.Range("J8").Value = 0.1846544489960017
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.987076
$ws.Range("N8").Value = 11.961228
$ws.Range("O8").Value = 0.2813308272685638
$ws.Range("P8").Value = 0.2813308272685638
$ws.Range("Q8").Value = 640.445943637708
$ws.Range("R8").Value = 5764.013492739372
$ws.Range("S8").Value = 0.05194898889486597
$ws.Range("T8").Value = 0.05194898889486597
$ws.Range("D9").Value = "MuSCs"
$ws.Range("G9").Value = 160.630483
$ws.Range("H9").Value = 481.891449
$ws.Range("I9").Value = 0.1846544489960017
$ws.Range("J9").Value = 0.1846544489960017
$ws.Range("M9").Value = 10.131229
$ws.Range("N9").Value = 30.393687
$ws.Range("O9").Value = 0.7148664925918803
$ws.Range("P9").Value = 0.7148664925918804
$ws.Range("Q9").Value = 1627.384207653607
$ws.Range("R9").Value = 14646.45786888246
$ws.Range("S9").Value = 0.132003278295258
$ws.Range("T9").Value = 0.132003278295258
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("G10").Value = 160.630483
$ws.Range("H10").Value = 481.891449
$ws.Range("I10").Value = 0.1846544489960017
$ws.Range("J10").Value = 0.1846544489960017
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.05389233333333333
$ws.Range("N10").Value = 0.161677
$ws.Range("O10").Value = 0.00380268013955587
$ws.Range("P10").Value = 0.00380268013955587
$ws.Range("Q10").Value = 8.656751533330333
$ws.Range("R10").Value = 77.91076379997298
$ws.Range("S10").Value = 0.000702181805877728
$ws.Range("T10").Value = 0.0007021818058777281
$ws.Range("D11").Value = "FAPs"
$ws.Range("G11").Value = 0.5357470000000001
$ws.Range("H11").Value = 1.607241
$ws.Range("I11").Value = 0.0006158735579862568
$ws.Range("J11").Value = 0.0006158735579862568
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.987076
$ws.Range("N11").Value = 11.961228
$ws.Range("O11").Value = 0.2813308272685638
$ws.Range("P11").Value = 0.2813308272685638
$ws.Range("Q11").Value = 2.136064005772
$ws.Range("R11").Value = 19.224576051948
$ws.Range("S11").Value = 0.0001732642175611074
$ws.Range("T11").Value = 0.0001732642175611074
$ws.Range("D12").Value = "MuSCs"
$ws.Range("G12").Value = 0.5357470000000001
$ws.Range("H12").Value = 1.607241
$ws.Range("I12").Value = 0.0006158735579862568
$ws.Range("J12").Value = 0.0006158735579862568
$ws.Range("M12").Value = 10.131229
$ws.Range("N12").Value = 30.393687
$ws.Range("O12").Value = 0.7148664925918803
$ws.Range("P12").Value = 0.7148664925918804
$ws.Range("Q12").Value = 5.427775543063
$ws.Range("R12").Value = 48.849979887567
$ws.Range("S12").Value = 0.0004402673702777174
$ws.Range("T12").Value = 0.0004402673702777174
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 0.5357470000000001
$ws.Range("H13").Value = 1.607241
$ws.Range("I13").Value = 0.0006158735579862568
$ws.Range("J13").Value = 0.0006158735579862568
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.05389233333333333
$ws.Range("N13").Value = 0.161677
$ws.Range("O13").Value = 0.00380268013955587
$ws.Range("P13").Value = 0.00380268013955587
$ws.Range("Q13").Value = 0.02887265590633334
$ws.Range("R13").Value = 0.259853903157
$ws.Range("S13").Value = 0.000002341970147431949
$ws.Range("T13").Value = 0.00000234197014743195
